$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B32").Value = 335
$ws.Range("B33").Value = 275
$ws.Range("B34").Value = 309
$ws.Range("B35").Value = 423
$ws.Range("B37").Value = 489
$ws.Range("B45").Value = 429
$ws.Range("B46").Value = 439
$ws.Range("B50").Value = 445

$ws.Range("A51").Value = 50
$ws.Range("B51").Value = 431
